$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings are not converted to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.121.83'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").Value = '1.670.56'
$ws.Range("E3").Value = '  -1.36%  '
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("D5").Value = '210.78'
$ws.Range("E5").Value = '  -3.31%  '
$ws.Range("D6").Value = '0.5255'
$ws.Range("E6").Value = '  -2.97%  '
$ws.Range("E7").Value = '  -0.61%  '
$ws.Range("E8").Value = '  -3.79%  '
$ws.Range("D9").Value = '0.06310'
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("D10").Value = '21.19'
$ws.Range("E10").Value = '  -2.24%  '
$ws.Range("D11").Value = '0.07560'
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '1.670.77'
$ws.Range("E12").Value = '  -1.73%  '
$ws.Range("D13").Value = '4.443'
$ws.Range("E13").Value = '  -2.07%  '
$ws.Range("D14").Value = '0.5576'
$ws.Range("E14").Value = '  -3.99%  '
$ws.Range("D15").Value = '67.03'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").Value = '0.000007955'
$ws.Range("E16").Value = '  -5.41%  '
$ws.Range("D17").Value = '26.163.44'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '4.750'
$ws.Range("E19").Value = '  -3.45%  '
$ws.Range("D20").Value = '187.01'
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("E21").Value = '  -4.63%  '
$ws.Range("D22").Value = '6.191'
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  -0.63%  '
$ws.Range("D24").Value = '149.58'
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("D26").Value = '7.519'
$ws.Range("D27").Value = '15.99'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '0.06277'
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("D29").Value = '1.358'
$ws.Range("E29").Value = '  -2.00%  '
$ws.Range("D30").Value = '1.282'
$ws.Range("E30").Value = '  -3.32%  '
$ws.Range("D31").Value = '3.512'
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("D32").Value = '3.419'
$ws.Range("E32").Value = '  -4.78%  '
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").Value = '  -3.14%  '
$ws.Range("D35").Value = '0.6057'
$ws.Range("E35").Value = '  -2.25%  '
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("D38").Value = '6.142'
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.101.91'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01612'
$ws.Range("E40").Value = '  -2.55%  '
$ws.Range("D41").Value = '0.8750'
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("D43").Value = '100.21'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").Value = '1.823.12'
$ws.Range("E44").Value = '  -1.23%  '
$ws.Range("D45").Value = '0.00000000109'
$ws.Range("E45").Value = '  -2.59%  '
$ws.Range("D46").Value = '55.50'
$ws.Range("D47").Value = '1.003'
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").Value = '8.038'
$ws.Range("E48").Value = '  -1.86%  '
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("E50").Value = '  -1.27%  '
$ws.Range("D51").Value = '5.989'
$ws.Range("E51").Value = '  -1.58%  '

# Restore default style on column D so no stray style attribute remains on cells
$ws.Range("D2:D51").Style = "Normal"

